$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 and 38: Dai and Fetch.AI swap positions, with updated prices
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'3.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.56%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.26%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D2").Value = "'71.027.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.83%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.648.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.59%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'605.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.51%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'198.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.38%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.26%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.07%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.219"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +8.74%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.648"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.10%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'53.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.42%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0000307"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.89%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'9.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.78%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.214.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +3.28%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'606.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.69%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'13.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.87%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'71.047.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.64%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.639.56"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +3.12%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'19.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.46%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.99%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  +0.96%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'18.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +1.13%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.90%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'103.96"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.63%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'4.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.59%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'3.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -5.14%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.15%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'9.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.73%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'33.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.14%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +12.49%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'7.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +2.84%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'12.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.89%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.95%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'63.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.54%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0₃0893"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +5.42%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'3.991.22"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +7.45%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D39").Value = "'519.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +6.93%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.392"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.14%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'36.83"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.73%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'3.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.94%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.05%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0463"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.92%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'3.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +5.41%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'2.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.15%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.62%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'8.62"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +2.04%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -0.41%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  +2.34%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.89%  "
$ws.Range("E51").Style = "Normal"
